$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2023-07-04 17:04:03"

$ws.Range("F2").Value = "OUT"
